$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 20.70050833333333
$ws.Range("H2").Value = 62.101525
$ws.Range("I2").Value = 0.8277101186170105
$ws.Range("J2").Value = 0.8277101186170105
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 68.00339
$ws.Range("N2").Value = 204.01017
$ws.Range("O2").Value = 0.6265962299909886
$ws.Range("P2").Value = 0.6265962299909885
$ws.Range("Q2").Value = 1407.704741389917
$ws.Range("R2").Value = 12669.34267250925
$ws.Range("S2").Value = 0.5186400398508128
$ws.Range("T2").Value = 0.5186400398508126

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 20.70050833333333
$ws.Range("H3").Value = 62.101525
$ws.Range("I3").Value = 0.8277101186170105
$ws.Range("J3").Value = 0.8277101186170105
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.99153
$ws.Range("N3").Value = 26.97459
$ws.Range("O3").Value = 0.08284967558015671
$ws.Range("P3").Value = 0.08284967558015671
$ws.Range("Q3").Value = 186.1292416944166
$ws.Range("R3").Value = 1675.16317524975
$ws.Range("S3").Value = 0.06857551480183235
$ws.Range("T3").Value = 0.06857551480183235

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 20.70050833333333
$ws.Range("H4").Value = 62.101525
$ws.Range("I4").Value = 0.8277101186170105
$ws.Range("J4").Value = 0.8277101186170105
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 10.76843933333333
$ws.Range("N4").Value = 32.305318
$ws.Range("O4").Value = 0.09922245772090688
$ws.Range("P4").Value = 0.09922245772090688
$ws.Range("Q4").Value = 222.9121681566611
$ws.Range("R4").Value = 2006.20951340995
$ws.Range("S4").Value = 0.08212743224964314
$ws.Range("T4").Value = 0.08212743224964314

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 20.70050833333333
$ws.Range("H5").Value = 62.101525
$ws.Range("I5").Value = 0.8277101186170105
$ws.Range("J5").Value = 0.8277101186170105
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 20.764887
$ws.Range("N5").Value = 62.294661
$ws.Range("O5").Value = 0.1913316367079478
$ws.Range("P5").Value = 0.1913316367079478
$ws.Range("Q5").Value = 429.8437163842249
$ws.Range("R5").Value = 3868.593447458024
$ws.Range("S5").Value = 0.1583671317147222
$ws.Range("T5").Value = 0.1583671317147222

$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.7925996666666667
$ws.Range("H6").Value = 2.377799
$ws.Range("I6").Value = 0.03169210888681734
$ws.Range("J6").Value = 0.03169210888681734
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 68.00339
$ws.Range("N6").Value = 204.01017
$ws.Range("O6").Value = 0.6265962299909886
$ws.Range("P6").Value = 0.6265962299909885
$ws.Range("Q6").Value = 53.89946424620333
$ws.Range("R6").Value = 485.0951782158299
$ws.Range("S6").Value = 0.01985815594894365
$ws.Range("T6").Value = 0.01985815594894365

$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.7925996666666667
$ws.Range("H7").Value = 2.377799
$ws.Range("I7").Value = 0.03169210888681734
$ws.Range("J7").Value = 0.03169210888681734
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 8.99153
$ws.Range("N7").Value = 26.97459
$ws.Range("O7").Value = 0.08284967558015671
$ws.Range("P7").Value = 0.08284967558015671
$ws.Range("Q7").Value = 7.126683680823333
$ws.Range("R7").Value = 64.14015312741
$ws.Range("S7").Value = 0.002625680939723818
$ws.Range("T7").Value = 0.002625680939723818

$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.7925996666666667
$ws.Range("H8").Value = 2.377799
$ws.Range("I8").Value = 0.03169210888681734
$ws.Range("J8").Value = 0.03169210888681734
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 10.76843933333333
$ws.Range("N8").Value = 32.305318
$ws.Range("O8").Value = 0.09922245772090688
$ws.Range("P8").Value = 0.09922245772090688
$ws.Range("Q8").Value = 8.535061426120222
$ws.Range("R8").Value = 76.815552835082
$ws.Range("S8").Value = 0.00314456893410861
$ws.Range("T8").Value = 0.00314456893410861

$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.7925996666666667
$ws.Range("H9").Value = 2.377799
$ws.Range("I9").Value = 0.03169210888681734
$ws.Range("J9").Value = 0.03169210888681734
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 20.764887
$ws.Range("N9").Value = 62.294661
$ws.Range("O9").Value = 0.1913316367079478
$ws.Range("P9").Value = 0.1913316367079478
$ws.Range("Q9").Value = 16.458242514571
$ws.Range("R9").Value = 148.124182631139
$ws.Range("S9").Value = 0.006063703064041259
$ws.Range("T9").Value = 0.006063703064041259

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.562510666666667
$ws.Range("H10").Value = 4.687532
$ws.Range("I10").Value = 0.06247701111592723
$ws.Range("J10").Value = 0.06247701111592723
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 68.00339
$ws.Range("N10").Value = 204.01017
$ws.Range("O10").Value = 0.6265962299909886
$ws.Range("P10").Value = 0.6265962299909885
$ws.Range("Q10").Value = 106.2560222444933
$ws.Range("R10").Value = 956.3042002004399
$ws.Range("S10").Value = 0.03914785962634509
$ws.Range("T10").Value = 0.03914785962634509

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.562510666666667
$ws.Range("H11").Value = 4.687532
$ws.Range("I11").Value = 0.06247701111592723
$ws.Range("J11").Value = 0.06247701111592723
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 8.99153
$ws.Range("N11").Value = 26.97459
$ws.Range("O11").Value = 0.08284967558015671
$ws.Range("P11").Value = 0.08284967558015671
$ws.Range("Q11").Value = 14.04936153465333
$ws.Range("R11").Value = 126.44425381188
$ws.Range("S11").Value = 0.005176200102172415
$ws.Range("T11").Value = 0.005176200102172416

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.562510666666667
$ws.Range("H12").Value = 4.687532
$ws.Range("I12").Value = 0.06247701111592723
$ws.Range("J12").Value = 0.06247701111592723
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 10.76843933333333
$ws.Range("N12").Value = 32.305318
$ws.Range("O12").Value = 0.09922245772090688
$ws.Range("P12").Value = 0.09922245772090688
$ws.Range("Q12").Value = 16.82580132168622
$ws.Range("R12").Value = 151.432211895176
$ws.Range("S12").Value = 0.006199122593978718
$ws.Range("T12").Value = 0.006199122593978719

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.562510666666667
$ws.Range("H13").Value = 4.687532
$ws.Range("I13").Value = 0.06247701111592723
$ws.Range("J13").Value = 0.06247701111592723
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 20.764887
$ws.Range("N13").Value = 62.294661
$ws.Range("O13").Value = 0.1913316367079478
$ws.Range("P13").Value = 0.1913316367079478
$ws.Range("Q13").Value = 32.445357429628
$ws.Range("R13").Value = 292.008216866652
$ws.Range("S13").Value = 0.011953828793431
$ws.Range("T13").Value = 0.011953828793431

$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1.953751
$ws.Range("H14").Value = 5.861253
$ws.Range("I14").Value = 0.07812076138024482
$ws.Range("J14").Value = 0.07812076138024483
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 68.00339
$ws.Range("N14").Value = 204.01017
$ws.Range("O14").Value = 0.6265962299909886
$ws.Range("P14").Value = 0.6265962299909885
$ws.Range("Q14").Value = 132.86169121589
$ws.Range("R14").Value = 1195.75522094301
$ws.Range("S14").Value = 0.04895017456488703
$ws.Range("T14").Value = 0.04895017456488703

$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 1.953751
$ws.Range("H15").Value = 5.861253
$ws.Range("I15").Value = 0.07812076138024482
$ws.Range("J15").Value = 0.07812076138024483
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 8.99153
$ws.Range("N15").Value = 26.97459
$ws.Range("O15").Value = 0.08284967558015671
$ws.Range("P15").Value = 0.08284967558015671
$ws.Range("Q15").Value = 17.56721072902999
$ws.Range("R15").Value = 158.10489656127
$ws.Range("S15").Value = 0.006472279736428119
$ws.Range("T15").Value = 0.00647227973642812

$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 1.953751
$ws.Range("H16").Value = 5.861253
$ws.Range("I16").Value = 0.07812076138024482
$ws.Range("J16").Value = 0.07812076138024483
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 10.76843933333333
$ws.Range("N16").Value = 32.305318
$ws.Range("O16").Value = 0.09922245772090688
$ws.Range("P16").Value = 0.09922245772090688
$ws.Range("Q16").Value = 21.03884911593933
$ws.Range("R16").Value = 189.349642043454
$ws.Range("S16").Value = 0.007751333943176396
$ws.Range("T16").Value = 0.007751333943176397

$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 1.953751
$ws.Range("H17").Value = 5.861253
$ws.Range("I17").Value = 0.07812076138024482
$ws.Range("J17").Value = 0.07812076138024483
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 20.764887
$ws.Range("N17").Value = 62.294661
$ws.Range("O17").Value = 0.1913316367079478
$ws.Range("P17").Value = 0.1913316367079478
$ws.Range("Q17").Value = 40.56941874113699
$ws.Range("R17").Value = 365.124768670233
$ws.Range("S17").Value = 0.01494697313575328
$ws.Range("T17").Value = 0.01494697313575328

